# Generate Report for Handback
# Adds a new handback entry (cc515795-00a0-4335-9ec6-ff4d30aa0e5e) as row 3 to all
# three worksheets (Overview, zh-cn, de-de), alongside updating the existing entry's
# file identifier from 3daeb6c3-5e79-4bff-a0e7-1d5c3314f394 to
# 59199e9b-216c-4568-9c85-c61bf9ca802a and refreshing the handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "3daeb6c3-5e79-4bff-a0e7-1d5c3314f394"
$newGuid = "59199e9b-216c-4568-9c85-c61bf9ca802a"
$addGuid = "cc515795-00a0-4335-9ec6-ff4d30aa0e5e"

$zhcnXliffHash = "bb4f61a0eeaf7b4a1715626c371b57790a83e8c1"
$dedeXliffHash = "bb4f61a0eeaf7b4a1715626c371b57790a83e8c1"
$addXliffHash  = "25ed7a593b7068aebc82aa7099dcafe3e64582a2"

$statusText = "Handed back: in sync with en-US"

# -----------------------------------------------------------------------
# Sheet 1: "Overview"
# -----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("G2").Value = "2016-08-20 23:05:11"

$wsOverview.Range("A3").Value = ($addGuid + ".md")
$wsOverview.Range("B3").Value = ("e2e\" + $addGuid + ".md")
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = "2016-08-20 23:05:11"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e98a1264ebdeb14f7d3c8bd3f1a05e36d8041b3e/e2e/$newGuid.md", $null, $null, ("e2e\" + $newGuid + ".md")) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e98a1264ebdeb14f7d3c8bd3f1a05e36d8041b3e/e2e/$addGuid.md", $null, $null, ("e2e\" + $addGuid + ".md")) | Out-Null

$wsOverview.Range("A1").Value = $wsOverview.Range("A1").Value

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3")) | Out-Null

# -----------------------------------------------------------------------
# Sheet 2: "zh-cn"
# -----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhTargetExisting = "$newGuid.$zhcnXliffHash.zh-cn.xlf"
$zhTargetNew      = "$addGuid.$addXliffHash.zh-cn.xlf"

# Update existing row (row 2)
$wsZhCn.Range("A2").Value = ($newGuid + ".md")
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("D2").Value = "e2e"
$wsZhCn.Range("E2").Value = "ht"
$wsZhCn.Range("F2").Value = "False"
$wsZhCn.Range("G2").Value = $zhTargetExisting
$wsZhCn.Range("H2").Value = "2016-08-20 23:05:06"
$wsZhCn.Range("I2").Value = ($newGuid + ".md")
$wsZhCn.Range("J2").Value = $zhTargetExisting
$wsZhCn.Range("K2").Value = "2016-08-20 23:05:30"
$wsZhCn.Range("M2").Value = "True"
$wsZhCn.Range("O2").Value = "False"

# New row (row 3)
$wsZhCn.Range("A3").Value = ($addGuid + ".md")
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = $zhTargetNew
$wsZhCn.Range("H3").Value = "2016-08-20 23:05:06"
$wsZhCn.Range("I3").Value = ($addGuid + ".md")
$wsZhCn.Range("J3").Value = $zhTargetNew
$wsZhCn.Range("K3").Value = "2016-08-20 23:05:30"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e98a1264ebdeb14f7d3c8bd3f1a05e36d8041b3e/e2e/$newGuid.md", $null, $null, ($newGuid + ".md")) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/94498bb75457e3feef569ccb419307102bda9641/e2e/$newGuid.md", $null, $null, ($newGuid + ".md")) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e98a1264ebdeb14f7d3c8bd3f1a05e36d8041b3e/e2e/$addGuid.md", $null, $null, ($addGuid + ".md")) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/94498bb75457e3feef569ccb419307102bda9641/e2e/$addGuid.md", $null, $null, ($addGuid + ".md")) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3")) | Out-Null

# -----------------------------------------------------------------------
# Sheet 3: "de-de"
# -----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deTargetExisting = "$newGuid.$dedeXliffHash.de-de.xlf"
$deTargetNew      = "$addGuid.$addXliffHash.de-de.xlf"

# Update existing row (row 2)
$wsDeDe.Range("A2").Value = ($newGuid + ".md")
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("D2").Value = "e2e"
$wsDeDe.Range("E2").Value = "ht"
$wsDeDe.Range("F2").Value = "False"
$wsDeDe.Range("G2").Value = $deTargetExisting
$wsDeDe.Range("H2").Value = "2016-08-20 23:05:11"
$wsDeDe.Range("I2").Value = ($newGuid + ".md")
$wsDeDe.Range("J2").Value = $deTargetExisting
$wsDeDe.Range("K2").Value = "2016-08-20 23:05:37"
$wsDeDe.Range("M2").Value = "True"
$wsDeDe.Range("O2").Value = "False"

# New row (row 3)
$wsDeDe.Range("A3").Value = ($addGuid + ".md")
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = $deTargetNew
$wsDeDe.Range("H3").Value = "2016-08-20 23:05:11"
$wsDeDe.Range("I3").Value = ($addGuid + ".md")
$wsDeDe.Range("J3").Value = $deTargetNew
$wsDeDe.Range("K3").Value = "2016-08-20 23:05:37"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e98a1264ebdeb14f7d3c8bd3f1a05e36d8041b3e/e2e/$newGuid.md", $null, $null, ($newGuid + ".md")) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/671bcb3bda9996bb7a274b833805e506d6d128a0/e2e/$newGuid.md", $null, $null, ($newGuid + ".md")) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e98a1264ebdeb14f7d3c8bd3f1a05e36d8041b3e/e2e/$addGuid.md", $null, $null, ($addGuid + ".md")) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/671bcb3bda9996bb7a274b833805e506d6d128a0/e2e/$addGuid.md", $null, $null, ($addGuid + ".md")) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3")) | Out-Null

Write-Host "Handback report rows generated."
